$d = $word.ActiveDocument

# The heading paragraph "3.1.3 日志配置" has a duplicated, trailing run
# containing the text "配置" immediately after the "_GoBack" bookmark
# (which must be kept in place). Remove just that trailing run's text
# so the heading reads "3.1.3 日志配置" instead of "3.1.3 日志配置配置".
$bm = $d.Bookmarks("_GoBack")
$pos = $bm.End
$extra = $d.Range($pos, $pos + 2)
$extra.Delete()
